$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "ATTENDANCE" column is inserted right before the existing "SCORE"
# column (old column F), pushing SCORE to column G. Raw attendance ("V")
# is filled in for every student row (6-10).
$ws.Range("F5").Value = "ATTENDANCE"
$ws.Range("G5").Value = "SCORE"
$ws.Range("F6:F10").Value = "V"
